# University workbook update
# - Notification system: clear out the old demo notification row and reset
#   the notification counter.
# - Rename the first teacher "Sam" -> "Matt" (and give him a last name,
#   "Peterson") and add a second teacher, "Sam Henchkins", together with a
#   second teacher login (lecturer2 / Lucy).
# - Misc bookkeeping (no. of teachers counters, a couple of leftover demo
#   rows in student_courses, a changed password owner in student_pswd).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# notifications sheet
# ---------------------------------------------------------------------
$notifications = $wb.Worksheets.Item("notifications")
$notifications.Activate()
$notifications.Range("A2:E2").ClearContents()
$notifications.Range("H3").Value = 0
$notifications.Range("F23").Select()

# ---------------------------------------------------------------------
# courses sheet
# The "Maths" course instructor is the teacher renamed below (Sam -> Matt),
# so the instructor text needs to follow along too.
# ---------------------------------------------------------------------
$courses = $wb.Worksheets.Item("courses")
$courses.Activate()
$courses.Range("E2").Value = "Matt"
$courses.Range("E2").Select()

# ---------------------------------------------------------------------
# student_courses sheet
# ---------------------------------------------------------------------
$studentCourses = $wb.Worksheets.Item("student_courses")
$studentCourses.Activate()
$studentCourses.Range("A2:E2").ClearContents()
$studentCourses.Range("A3:E3").ClearContents()
$studentCourses.Range("A4:E4").ClearContents()
$studentCourses.Range("L4").Value = 0
$studentCourses.Range("L4").Select()

# ---------------------------------------------------------------------
# student_pswd sheet
# ---------------------------------------------------------------------
$studentPswd = $wb.Worksheets.Item("student_pswd")
$studentPswd.Activate()
$studentPswd.Range("C3").Value = "Dave"
$studentPswd.Range("G6").Select()

# ---------------------------------------------------------------------
# Teachers sheet
# ---------------------------------------------------------------------
$teachers = $wb.Worksheets.Item("Teachers")
$teachers.Activate()
$teachers.Range("B2").Value = "Matt"
$teachers.Range("C2").Value = "Peterson"
$teachers.Range("A3").Value = 2
$teachers.Range("B3").Value = "Sam"
$teachers.Range("C3").Value = "Henchkins"
$teachers.Range("D3").Value = "h"
$teachers.Range("E3").Value = "jh"
$teachers.Range("H4").Value = 2
$teachers.Range("C5").Select()

# ---------------------------------------------------------------------
# teacher_psswd sheet
# ---------------------------------------------------------------------
$teacherPsswd = $wb.Worksheets.Item("teacher_psswd")
$teacherPsswd.Activate()
$teacherPsswd.Range("C2").Value = "Matt"
$teacherPsswd.Range("A3").Value = 2
$teacherPsswd.Range("B3").Value = "lecturer2"
$teacherPsswd.Range("C3").Value = "Lucy"
$teacherPsswd.Range("F4").Value = 2

# ---------------------------------------------------------------------
# Leave the workbook focused back on the notifications sheet, which was
# the active tab before editing.
# ---------------------------------------------------------------------
$notifications.Activate()
$notifications.Range("F23").Select()
